$wb = $excel.ActiveWorkbook

# Update "展览" sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5278
$ws1.Range("F3").Value = 163
$ws1.Range("F4").Value = 917

# Update "全部类型" sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5278
$ws4.Range("F3").Value = 163
$ws4.Range("F4").Value = 917
